$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 97.26667
$ws.Range("I33").Value = 99.92856999999999
$ws.Range("J33").Value = 60
$ws.Range("K33").Value = 99.92856999999999
$ws.Range("L33").Value = 60
$ws.Range("M33").Value = 129.07143
$ws.Range("N33").Value = -518

$ws.Range("H76").Value = 3192.8965
$ws.Range("I76").Value = 3220.8333
$ws.Range("J76").Value = 3058.8
$ws.Range("K76").Value = 3220.8333
$ws.Range("L76").Value = 3058.8
$ws.Range("M76").Value = -2905.8333
$ws.Range("N76").Value = -3688.8

$ws.Range("H79").Value = 3192.8965
$ws.Range("I79").Value = 3220.8333
$ws.Range("J79").Value = 3058.8
$ws.Range("K79").Value = 3220.8333
$ws.Range("L79").Value = 3058.8
$ws.Range("M79").Value = -2128.8333
$ws.Range("N79").Value = -5242.8

$ws.Range("H129").Value = 1333.772
$ws.Range("I129").Value = 438.61905
$ws.Range("J129").Value = 1855.9445
$ws.Range("K129").Value = 1315.85715
$ws.Range("L129").Value = 5567.833500000001
$ws.Range("M129").Value = 3684.14285
$ws.Range("N129").Value = -15567.8335

$ws.Range("H131").Value = 4482.643
$ws.Range("I131").Value = 492.5
$ws.Range("J131").Value = 5147.6665
$ws.Range("K131").Value = 1477.5
$ws.Range("L131").Value = 15442.9995
$ws.Range("M131").Value = 3562.5
$ws.Range("N131").Value = -25522.9995

$ws.Range("H132").Value = 2534.6938
$ws.Range("I132").Value = 2431.8408
$ws.Range("J132").Value = 3439.8
$ws.Range("K132").Value = 7295.5224
$ws.Range("L132").Value = 10319.4
$ws.Range("M132").Value = -4765.5224
$ws.Range("N132").Value = -15379.4

$ws.Range("H137").Value = 1232.151
$ws.Range("I137").Value = 1124.6666
$ws.Range("J137").Value = 1409.5
$ws.Range("K137").Value = 3373.9998
$ws.Range("L137").Value = 4228.5
$ws.Range("M137").Value = -823.9998000000001
$ws.Range("N137").Value = -9328.5

$ws.Range("H138").Value = 1466.85
$ws.Range("I138").Value = 637.7442
$ws.Range("J138").Value = 2092.3157
$ws.Range("K138").Value = 1913.2326
$ws.Range("L138").Value = 6276.9471
$ws.Range("M138").Value = 3226.7674
$ws.Range("N138").Value = -16556.9471


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2411.3572
$ws.Range("I61").Value = 2046.1025
$ws.Range("J61").Value = 3249.2942
$ws.Range("K61").Value = 2046.1025
$ws.Range("L61").Value = 3249.2942
$ws.Range("M61").Value = -1834.1025
$ws.Range("N61").Value = -3673.2942

$ws.Range("H63").Value = 24701.477
$ws.Range("I63").Value = 78140.75
$ws.Range("J63").Value = 4661.75
$ws.Range("K63").Value = 78140.75
$ws.Range("L63").Value = 4661.75
$ws.Range("M63").Value = -77454.75
$ws.Range("N63").Value = -6033.75

$ws.Range("H66").Value = 24701.477
$ws.Range("I66").Value = 78140.75
$ws.Range("J66").Value = 4661.75
$ws.Range("K66").Value = 390703.75
$ws.Range("L66").Value = 23308.75
$ws.Range("M66").Value = -387271.75
$ws.Range("N66").Value = -30172.75

$ws.Range("H136").Value = 2411.3572
$ws.Range("I136").Value = 2046.1025
$ws.Range("J136").Value = 3249.2942
$ws.Range("K136").Value = 6138.3075
$ws.Range("L136").Value = 9747.882599999999
$ws.Range("M136").Value = -3588.3075
$ws.Range("N136").Value = -14847.8826


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2096.8
$ws.Range("I20").Value = 1732.8667
$ws.Range("J20").Value = 2642.7
$ws.Range("K20").Value = 1732.8667
$ws.Range("L20").Value = 2642.7
$ws.Range("M20").Value = -1485.8667
$ws.Range("N20").Value = -3136.7

$ws.Range("H134").Value = 2641.0454
$ws.Range("I134").Value = 2328.2856
$ws.Range("J134").Value = 3188.375
$ws.Range("K134").Value = 6984.8568
$ws.Range("L134").Value = 9565.125
$ws.Range("M134").Value = -4449.8568
$ws.Range("N134").Value = -14635.125


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 29250.5
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 29250.5
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 29250.5
$ws.Range("N4").Value = -29474.5
$ws.Range("M4").ClearContents()

$ws.Range("H31").Value = 4324.314
$ws.Range("I31").Value = 1211.84
$ws.Range("J31").Value = 7317.077
$ws.Range("K31").Value = 1211.84
$ws.Range("L31").Value = 7317.077
$ws.Range("M31").Value = -916.8399999999999
$ws.Range("N31").Value = -7907.077

$ws.Range("H34").Value = 4324.314
$ws.Range("I34").Value = 1211.84
$ws.Range("J34").Value = 7317.077
$ws.Range("K34").Value = 1211.84
$ws.Range("L34").Value = 7317.077
$ws.Range("M34").Value = -1009.84
$ws.Range("N34").Value = -7721.077

$ws.Range("H58").Value = 1258.0541
$ws.Range("I58").Value = 1045.1666
$ws.Range("J58").Value = 1459.7368
$ws.Range("K58").Value = 1045.1666
$ws.Range("L58").Value = 1459.7368
$ws.Range("M58").Value = -842.1666
$ws.Range("N58").Value = -1865.7368

$ws.Range("H99").Value = 1884.3334

$ws.Range("H126").Value = 1884.3334

$ws.Range("H132").Value = 9806483
$ws.Range("I132").Value = 2752.5715
$ws.Range("J132").Value = 55557224
$ws.Range("K132").Value = 8257.7145
$ws.Range("L132").Value = 166671672
$ws.Range("M132").Value = -5727.7145
$ws.Range("N132").Value = -166676732

$ws.Range("H134").Value = 3740.4187
$ws.Range("I134").Value = 4109.543
$ws.Range("J134").Value = 2125.5
$ws.Range("K134").Value = 12328.629
$ws.Range("L134").Value = 6376.5
$ws.Range("M134").Value = -9793.628999999999

$ws.Range("H136").Value = 1258.0541
$ws.Range("I136").Value = 1045.1666
$ws.Range("J136").Value = 1459.7368
$ws.Range("K136").Value = 3135.4998
$ws.Range("L136").Value = 4379.2104
$ws.Range("M136").Value = -585.4998000000001
$ws.Range("N136").Value = -9479.2104


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2835.0308
$ws.Range("I131").Value = 521.375
$ws.Range("J131").Value = 3159.7544
$ws.Range("K131").Value = 1564.125
$ws.Range("L131").Value = 9479.263199999999
$ws.Range("M131").Value = 3475.875
$ws.Range("N131").Value = -19559.2632

$ws.Range("H139").Value = 3316.85
$ws.Range("I139").Value = 1382.35
$ws.Range("J139").Value = 5251.35
$ws.Range("K139").Value = 4147.049999999999
$ws.Range("L139").Value = 15754.05
$ws.Range("M139").Value = 992.9500000000007
$ws.Range("N139").Value = -26034.05


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 65
$ws.Range("I2").Value = 65
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 65
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 48
$ws.Range("N2").ClearContents()

$ws.Range("H102").Value = 1741
$ws.Range("I102").Value = 1715.15
$ws.Range("J102").Value = 1999.5
$ws.Range("K102").Value = 1715.15
$ws.Range("L102").Value = 1999.5
$ws.Range("M102").Value = -93.15000000000009
$ws.Range("N102").Value = -5243.5

$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 6000
$ws.Range("N126").Value = -10940
$ws.Range("M126").ClearContents()

$ws.Range("H132").Value = 3834.889
$ws.Range("I132").Value = 5700
$ws.Range("J132").Value = 3302
$ws.Range("K132").Value = 17100
$ws.Range("L132").Value = 9906
$ws.Range("M132").Value = -14570
$ws.Range("N132").Value = -14966


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 3052502
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 3052502
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 3052502
$ws.Range("N2").Value = -3052726

$ws.Range("H7").Value = 4392.3687
$ws.Range("I7").Value = 4237.273
$ws.Range("J7").Value = 4605.625
$ws.Range("K7").Value = 4237.273
$ws.Range("L7").Value = 4605.625
$ws.Range("M7").Value = -4125.273
$ws.Range("N7").Value = -4829.625

$ws.Range("H40").Value = 335334.66
$ws.Range("I40").Value = 1000004
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 1000004
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -999868

$ws.Range("H68").Value = 1480.6786
$ws.Range("I68").Value = 1450.3914
$ws.Range("J68").Value = 1620
$ws.Range("K68").Value = 1450.3914
$ws.Range("L68").Value = 1620
$ws.Range("M68").Value = -701.3914
$ws.Range("N68").Value = -3118

$ws.Range("H71").Value = 1480.6786
$ws.Range("I71").Value = 1450.3914
$ws.Range("J71").Value = 1620
$ws.Range("K71").Value = 7251.957
$ws.Range("L71").Value = 8100
$ws.Range("M71").Value = -3507.957
$ws.Range("N71").Value = -15588

$ws.Range("H122").Value = 3863.3845
$ws.Range("I122").Value = 1804
$ws.Range("J122").Value = 4035
$ws.Range("K122").Value = 5412
$ws.Range("L122").Value = 12105
$ws.Range("M122").Value = -2962
$ws.Range("N122").Value = -17005

$ws.Range("H126").Value = 4392.3687
$ws.Range("I126").Value = 4237.273
$ws.Range("J126").Value = 4605.625
$ws.Range("K126").Value = 12711.819
$ws.Range("L126").Value = 13816.875
$ws.Range("M126").Value = -10241.819
$ws.Range("N126").Value = -18756.875

$ws.Range("H127").Value = 30749.17
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 30749.17
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 30749.17
$ws.Range("N127").Value = -40669.17

$ws.Range("H132").Value = 2818.3022
$ws.Range("I132").Value = 2795.7856
$ws.Range("J132").Value = 2860.3333
$ws.Range("K132").Value = 8387.356800000001
$ws.Range("L132").Value = 8580.999899999999
$ws.Range("M132").Value = -5857.356800000001
$ws.Range("N132").Value = -13640.9999

$ws.Range("H136").Value = 3473439.8
$ws.Range("I136").Value = 1167.675
$ws.Range("J136").Value = 20834800
$ws.Range("K136").Value = 3503.025
$ws.Range("L136").Value = 62504400
$ws.Range("M136").Value = -953.0249999999996
$ws.Range("N136").Value = -62509500


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2661.9048
$ws.Range("I122").Value = 2629.75
$ws.Range("J122").Value = 3305
$ws.Range("K122").Value = 7889.25
$ws.Range("L122").Value = 9915
$ws.Range("M122").Value = -5439.25

$ws.Range("H132").Value = 6412180
$ws.Range("I132").Value = 2604.7856
$ws.Range("J132").Value = 13890018
$ws.Range("K132").Value = 7814.3568
$ws.Range("L132").Value = 41670054
$ws.Range("M132").Value = -5284.3568
$ws.Range("N132").Value = -41675114

$ws.Range("H136").Value = 2062.2297
$ws.Range("I136").Value = 1842.8103
$ws.Range("J136").Value = 2857.625
$ws.Range("K136").Value = 5528.4309
$ws.Range("L136").Value = 8572.875
$ws.Range("M136").Value = -2978.4309
$ws.Range("N136").Value = -13672.875

